$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J3").Value = "2019.05.20"
$ws.Range("B5").Value = "安装配置mysql"
$ws.Range("B6").Value = "学习使用，代码的了解"
$ws.Range("B7").Value = "库和表的一些增删改查操作"

$ws.Range("J3:L3").Select()
